# Update "想去人数" (interested-count) values for several events on the
# "展览" and "全部类型" worksheets, per the regenerated gh-pages data dump.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5714
$ws1.Range("F5").Value = 963
$ws1.Range("F7").Value = 2621
$ws1.Range("F11").Value = 99
$ws1.Range("F13").Value = 2457
$ws1.Range("F14").Value = 504

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5714
$ws4.Range("F6").Value = 963
$ws4.Range("F8").Value = 2621
$ws4.Range("F13").Value = 99
$ws4.Range("F15").Value = 2457
$ws4.Range("F16").Value = 504

$wb.Save()
